# Add a new "plus one" problem entry to the "数组" (Array) worksheet (5th sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Column A: No.
$ws.Cells.Item(6, 1).Value2 = 5
# Column B: leetcode id
$ws.Cells.Item(6, 2).Value2 = 66

# Column C: 题目 (problem statement)
$problemText = @"
给定一个由整数组成的非空数组所表示的非负整数，在该数的基础上加一。 
 最高位数字存放在数组的首位， 数组中每个元素只存储单个数字。 
 你可以假设除了整数 0 之外，这个整数不会以零开头。 
 示例 1: 
 输入: [1,2,3]
输出: [1,2,4]
解释: 输入数组表示数字 123。
 示例 2: 
 输入: [4,3,2,1]
输出: [4,3,2,2]
解释: 输入数组表示数字 4321。
 Related Topics 数组
"@
$ws.Cells.Item(6, 3).Value2 = $problemText

# Column D: 解题方法 (solution approach)
$solutionText = @"
0 迭代数组
1 获取数值并且加一
     如果之和等于10，说明这个原数字是9，应该将这一位变成，放心：下一次的循环会直接digits[i] = digits[i] + 1;
     如果之和不扽与10，说明没有进位，就可以直接返回当前的结果
2 继续迭代
3 数组迭代完成，说明最高位还有进位，创建新数组保存
"@
$ws.Cells.Item(6, 4).Value2 = $solutionText

# Column E: 解题关键词 (keywords)
$ws.Cells.Item(6, 5).Value2 = "相加`n进位"

# Column F: 时间复杂度 (time complexity)
$ws.Cells.Item(6, 6).Value2 = "O(N)"

# Column G: 空间复杂度 (space complexity)
$ws.Cells.Item(6, 7).Value2 = "O(1)"

# Match the tall row height used by the other multi-line rows in this sheet.
$ws.Rows.Item(6).RowHeight = 409.6

# Update the view so the new row/cell is the active selection, like in the saved file.
$ws.Activate()
$ws.Range("E6").Select()
